# Daily attendance processing - 2025-11-17 03:03:09
# Rotate the "Recorded By" (column G) list so the last contributor listed
# moves to the front, for every data row on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $value = $cell.Value2

    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ", "
        if ($parts.Count -gt 1) {
            $rotated = @($parts[-1]) + $parts[0..($parts.Count - 2)]
            $cell.Value2 = [string]::Join(", ", $rotated)
        }
    }
}
